# Auto-generated edit script: updates country stats and reorders a few
# countries whose totals changed rank, matching the "Update countries &
# provincias Spain" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 16 de Julio de 2020 a las 17:06"
$ws.Cells.Item(4, 2).Value = 3634420
$ws.Cells.Item(4, 3).Value = 17593
$ws.Cells.Item(4, 4).Value = 1646698
$ws.Cells.Item(4, 5).Value = 1847356
$ws.Cells.Item(4, 7).Value = 222
$ws.Cells.Item(4, 8).Value = 140366
$ws.Cells.Item(6, 2).Value = 980710
$ws.Cells.Item(6, 3).Value = 10541
$ws.Cells.Item(6, 4).Value = 620194
$ws.Cells.Item(6, 5).Value = 335394
$ws.Cells.Item(6, 7).Value = 193
$ws.Cells.Item(6, 8).Value = 25122
$ws.Cells.Item(13, 2).Value = 292552
$ws.Cells.Item(13, 3).Value = 641
$ws.Cells.Item(13, 7).Value = 66
$ws.Cells.Item(13, 8).Value = 45119
$ws.Cells.Item(45, 4).Value = 43256
$ws.Cells.Item(45, 5).Value = 3843
$ws.Cells.Item(61, 2).Value = 20264
$ws.Cells.Item(61, 3).Value = 224
$ws.Cells.Item(61, 5).Value = 5958
$ws.Cells.Item(61, 7).Value = 7
$ws.Cells.Item(61, 8).Value = 666
$ws.Cells.Item(81, 2).Value = 8623
$ws.Cells.Item(81, 3).Value = 93
$ws.Cells.Item(81, 5).Value = 3657
$ws.Cells.Item(81, 7).Value = 8
$ws.Cells.Item(81, 8).Value = 401
$ws.Cells.Item(91, 2).Value = 6393
$ws.Cells.Item(91, 3).Value = 94
$ws.Cells.Item(91, 4).Value = 3818
$ws.Cells.Item(91, 5).Value = 2541
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 34
$ws.Cells.Item(108, 5).Value = 495
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 15
$ws.Cells.Item(110, 2).Value = 2686
$ws.Cells.Item(110, 3).Value = 15
$ws.Cells.Item(110, 5).Value = 668
$ws.Cells.Item(120, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(120, 3).Value = 60
$ws.Cells.Item(120, 4).Value = 773
$ws.Cells.Item(120, 5).Value = 1103
$ws.Cells.Item(120, 8).Value = 26
$ws.Cells.Item(121, 1).Value = "Lituania"
$ws.Cells.Item(121, 2).Value = 1902
$ws.Cells.Item(121, 3).Value = 6
$ws.Cells.Item(121, 4).Value = 1593
$ws.Cells.Item(121, 5).Value = 230
$ws.Cells.Item(121, 8).Value = 79
$ws.Cells.Item(122, 1).Value = "Eslovenia"
$ws.Cells.Item(122, 2).Value = 1897
$ws.Cells.Item(122, 3).Value = 19
$ws.Cells.Item(122, 4).Value = 1522
$ws.Cells.Item(122, 5).Value = 264
$ws.Cells.Item(122, 8).Value = 111
$ws.Cells.Item(123, 1).Value = "Zambia"
$ws.Cells.Item(123, 2).Value = 1895
$ws.Cells.Item(123, 4).Value = 1412
$ws.Cells.Item(123, 5).Value = 441
$ws.Cells.Item(123, 8).Value = 42
$ws.Cells.Item(132, 1).Value = "Mozambique"
$ws.Cells.Item(132, 2).Value = 1383
$ws.Cells.Item(132, 3).Value = 53
$ws.Cells.Item(132, 4).Value = 375
$ws.Cells.Item(132, 5).Value = 999
$ws.Cells.Item(132, 8).Value = 9
$ws.Cells.Item(133, 1).Value = "Benin"
$ws.Cells.Item(133, 2).Value = 1378
$ws.Cells.Item(133, 4).Value = 557
$ws.Cells.Item(133, 5).Value = 795
$ws.Cells.Item(133, 8).Value = 26
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"

Write-Output "applied 76 cell updates"
